# Updates the "cryptos" worksheet with refreshed Price / Volume(1h) figures
# (and, for a few re-ranked coins, refreshed Coin name + Link as well),
# matching the GitHub Actions data refresh described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark Price (D) cells that are being updated as Text format first so
# values such as "1.00" / "507.00" / "69.918.81" are preserved verbatim
# as text instead of being auto-coerced into numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.918.81"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "3.752.54"
$ws.Range("E3").Value = "  +2.61%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "622.42"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("D6").Value = "178.99"
$ws.Range("E6").Value = "  -1.46%  "
$ws.Range("D7").Value = "3.751.34"
$ws.Range("E7").Value = "  +2.66%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  -0.25%  "
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  +3.07%  "
$ws.Range("D11").Value = "6.32"
$ws.Range("E11").Value = "  -5.87%  "
$ws.Range("D12").Value = "0.489"
$ws.Range("E12").Value = "  -2.60%  "
$ws.Range("D13").Value = "41.11"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").Value = "0.0000260"
$ws.Range("E14").Value = "  +2.35%  "
$ws.Range("D15").Value = "4.371.11"
$ws.Range("E15").Value = "  +2.38%  "
$ws.Range("D16").Value = "3.750.16"
$ws.Range("E16").Value = "  +2.49%  "
$ws.Range("D17").Value = "69.945.96"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "0.124"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "7.63"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("D20").Value = "16.69"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "507.00"
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").Value = "9.44"
$ws.Range("E22").Value = "  +1.63%  "
$ws.Range("D23").Value = "0.729"
$ws.Range("E23").Value = "  -2.10%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("D25").Value = "87.18"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").Value = "13.18"
$ws.Range("E26").Value = "  -2.45%  "
$ws.Range("D27").Value = "11.13"
$ws.Range("E27").Value = "  +1.13%  "
$ws.Range("D28").Value = "0.0000136"
$ws.Range("E28").Value = "  +23.44%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").Value = "2.50"
$ws.Range("E30").Value = "  -2.16%  "
$ws.Range("E31").Value = "  -0.32%  "
$ws.Range("D32").Value = "7.92"
$ws.Range("E32").Value = "  -2.09%  "
$ws.Range("D33").Value = "31.61"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").Value = "0.115"
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("E35").Value = "  -0.15%  "
$ws.Range("E36").Value = "  +4.11%  "
$ws.Range("D37").Value = "6.22"
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").Value = "0.335"
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("D39").Value = "0.133"
$ws.Range("E39").Value = "  +1.76%  "
$ws.Range("E40").Value = "  -4.13%  "
$ws.Range("D41").Value = "50.21"
$ws.Range("E41").Value = "  -3.15%  "
$ws.Range("D42").Value = "45.22"
$ws.Range("E42").Value = "  -1.72%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "423.62"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D44").Value = "8.72"
$ws.Range("E44").Value = "  -1.41%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "2.85"
$ws.Range("E45").Value = "  +2.48%  "
$ws.Range("D46").Value = "3.009.94"
$ws.Range("E46").Value = "  -3.56%  "
$ws.Range("E47").Value = "  -1.98%  "
$ws.Range("D48").Value = "27.31"
$ws.Range("E48").Value = "  -4.19%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "138.46"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("D51").Value = "2.51"
$ws.Range("E51").Value = "  +1.78%  "
